# Timetable re-shuffle: move several course codes to their new
# room/time slots across the mon/tue/wed/thur/fri sheets.
# Each course normally occupies two adjacent 1-hour columns in a given
# room row, so cells are cleared at the old slot and written at the new
# slot (sometimes within the same row, sometimes a different room/day).

$wb = $excel.ActiveWorkbook

# ----- mon -----
$ws = $wb.Worksheets.Item("mon")
$ws.Range("I4").Value  = ""          # CSC423 removed from Biology lab 4-5pm
$ws.Range("B3").Value  = "CSC423"    #   -> now Bio Chem Lab 2, 8-9am
$ws.Range("F7").Value  = ""          # CSC111 removed from Chem Prac. Lab 12-1pm
$ws.Range("G7").Value  = ""
$ws.Range("D13").Value = "CSC424"    # CSC424 added to Hall 201, 10-12
$ws.Range("E13").Value = "CSC424"
$ws.Range("J14").Value = "CSC424"    # CSC424 added to Hall 202, 5-6pm
$ws.Range("J15").Value = ""          # BIO111 removed from Hall 203, 5-6pm
$ws.Range("D17").Value = ""          # CSC424 removed from Hall 306, 10-11am
$ws.Range("G17").Value = "BIO111"    #   BIO111 added to Hall 306, 1-2pm
$ws.Range("H19").Value = ""          # CSC424 removed from Hall 308, 3-5pm
$ws.Range("I19").Value = ""
$ws.Range("E20").Value = ""          # MAT111 removed from LT 1, 11-1pm
$ws.Range("F20").Value = ""
$ws.Range("F21").Value = ""          # GST111 removed from LT 2, 12-1pm
$ws.Range("G21").Value = "MAT111"    #   LT 2, 1-2pm GST111 -> MAT111

# ----- tue -----
$ws = $wb.Worksheets.Item("tue")
$ws.Range("B7").Value  = "CSC425"    # CSC425 added to Chem Prac. Lab, 8-10am
$ws.Range("C7").Value  = "CSC425"
$ws.Range("C14").Value = "CSC442"    # CSC442 added to Hall 202, 9-11am
$ws.Range("D14").Value = "CSC442"
$ws.Range("F20").Value = "EDS421"    # EDS421 added to LT 1, 12-1pm
$ws.Range("K20").Value = ""          # MAT111 removed from LT 1, 6-7pm
$ws.Range("C21").Value = "MAT112"    # MAT112 now also at LT 2, 9-10am
$ws.Range("E21").Value = ""          # MAT112 removed from LT 2, 11am-12pm
$ws.Range("I21").Value = ""          # CST111 removed from LT 2, 4-6pm
$ws.Range("J21").Value = ""

# ----- wed -----
$ws = $wb.Worksheets.Item("wed")
$ws.Range("I7").Value  = ""          # CSC425 removed from Chem Prac. Lab, 4-5pm
$ws.Range("D14").Value = ""          # BIO111 removed from Hall 202, 10-12
$ws.Range("E14").Value = ""
$ws.Range("F15").Value = "CSC424"    # CSC424 added to Hall 203, 12-1pm
$ws.Range("I17").Value = "CIS421"    # CIS421 added to Hall 306, 4-6pm
$ws.Range("J17").Value = "CIS421"
$ws.Range("D20").Value = ""          # CIT111 removed from LT 1, 10-12
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = "MAT111"    #   LT 1, 12-1pm & 1-2pm -> MAT111
$ws.Range("G20").Value = "MAT111"
$ws.Range("I20").Value = ""          # EDS421 removed from LT 1, 4-5pm
$ws.Range("J20").Value = "TMC421"    #   LT 1, 5-6pm -> TMC421
$ws.Range("F22").Value = "CSC441"    # CSC441 added to Micro Bio. Lab, 12-1pm
$ws.Range("G22").Value = "CSC441"
$ws.Range("H25").Value = "CSC111"    # CSC111 added to Physics Lab (ground floor), 3-5pm
$ws.Range("I25").Value = "CSC111"

# ----- thur -----
$ws = $wb.Worksheets.Item("thur")
$ws.Range("G2").Value  = ""          # CSC425 removed from Bio Chem Lab 1, 1-2pm & 3-4pm
$ws.Range("H2").Value  = ""
$ws.Range("F4").Value  = ""          # CSC111 removed from Biology lab, 12-1pm
$ws.Range("E13").Value = ""          # CSC424 removed from Hall 201, 11am-12pm
$ws.Range("D17").Value = ""          # CSC424 removed from Hall 306, 10-11am
$ws.Range("E17").Value = ""          #   and 11am-12pm
$ws.Range("F17").Value = ""          # CSC442 removed from Hall 306, 12-1pm
$ws.Range("F19").Value = "CSC424"    # CSC424 added to Hall 308, 12-1pm
$ws.Range("G19").Value = "CSC424"    #   and 1-2pm
$ws.Range("H20").Value = "CST111"    # CST111 added to LT 1, 3-4pm
$ws.Range("I20").Value = "CST111"    #   and 4-5pm
$ws.Range("G21").Value = "CSC111"    # CSC111 added to LT 2, 1-2pm
$ws.Range("E24").Value = ""          # CSC423 removed from Physics Lab (2nd floor), 11am-12pm
$ws.Range("F24").Value = ""          #   and 12-1pm

# ----- fri -----
$ws = $wb.Worksheets.Item("fri")
$ws.Range("F7").Value  = ""          # CIS421 removed from Chem Prac. Lab, 12-1pm
$ws.Range("G7").Value  = ""
$ws.Range("E10").Value = "CSC423"    # CSC423 added to Hall 107, 11am-1pm
$ws.Range("F10").Value = "CSC423"
$ws.Range("D13").Value = "CSC442"    # CSC442 added to Hall 201, 10-11am
$ws.Range("C20").Value = "TMC111"    # TMC111 added to LT 1, 9-10am
$ws.Range("D21").Value = "CHM111"    # CHM111 added to LT 2, 10-12
$ws.Range("E21").Value = "CHM111"
$ws.Range("C28").Value = ""          # CSC441 removed from Studio 3, 9-11am
$ws.Range("D28").Value = ""
